# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-38, columns E & F) is
# rebuilt in the opposite chronological order: what used to start at
# 1809 (Sep-2018) and end at 2007 (Jul-2020) now starts at 2007 and ends
# at 1809. The associated "Valor Mora" column is reversed the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2007","2006","2005","2004","2003","2002","2001","1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901","1812","1811","1810","1809")
$valores = @(22916,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
